$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update data values per diff
$ws2.Range("B3").Value = 95
$ws2.Range("C3").Value = 265
$ws2.Range("B7").Value = 93
$ws2.Range("B13").Value = 352
$ws2.Range("B17").Value = 324
$ws2.Range("J4").Value = 0.1444

# Title rows: merge A:D and center-align, matching Sheet1's header pattern
$ws2.Range("A1:D1").Merge()
$ws2.Range("A1:D1").HorizontalAlignment = -4108

$ws2.Range("A5:D5").Merge()
$ws2.Range("A5:D5").HorizontalAlignment = -4108

$ws2.Range("A11:D11").Merge()
$ws2.Range("A11:D11").HorizontalAlignment = -4108

$ws2.Range("A15:D15").Merge()
$ws2.Range("A15:D15").HorizontalAlignment = -4108
